# Add a "Tag List" column to the FundFormula sheet, between "Rule For" (D)
# and "Commitment Type" (old E, now F), and populate it with tag values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FundFormula")

# Insert a new column before the current column E ("Commitment Type"),
# shifting everything from E onward one column to the right.
$ws.Range("E1").EntireColumn.Insert()

# New column header + data values.
$ws.Range("E1").Value = "Tag List"
$ws.Range("E2").Value = "Monthly"
$ws.Range("E3").Value = "Monthly"
$ws.Range("E4").Value = "Quarterly"
$ws.Range("E5").Value = "Monthly"
$ws.Range("E6").Value = "Monthly"
$ws.Range("E7").Value = "Monthly"
$ws.Range("E8").Value = "Monthly"
$ws.Range("E9").Value = "Monthly"
$ws.Range("E10").Value = "Monthly"

# Match the column width used by the adjacent "Rule For" column rather than
# autofitting, since this column was typed in manually (not bestFit).
$ws.Range("E1").ColumnWidth = $ws.Range("D1").ColumnWidth

# Update the active selection to reflect the newly filled column.
$ws.Range("E5:E10").Select()

$wb.Save()
